# Update "Latest HO Xliff Generate Date" / "Correspond Handoff/Handback Datetime"
# timestamps to reflect the newly-generated handback report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview!G4 - Latest HO Xliff Generate Date (shared with de-de!H4)
$wsOverview.Range("G4").Value = "2016-10-14 07:43:30"

# zh-cn!H4 - Correspond Handoff Datetime
$wsZhCn.Range("H4").Value = "2016-10-14 07:43:19"

# zh-cn!K4 - Correspond Handback DateTime
$wsZhCn.Range("K4").Value = "2016-10-14 07:44:03"

# de-de!H4 - Correspond Handoff Datetime (shares string with Overview!G4)
$wsDeDe.Range("H4").Value = "2016-10-14 07:43:30"

# de-de!K4 - Correspond Handback DateTime
$wsDeDe.Range("K4").Value = "2016-10-14 07:44:20"
